# "Rendu final de la bataille navale version 1.0" — log the final
# delivery entry in the "Journal de Bord" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new journal entry (date + description), reusing the same
#     look (border/font/wrap/date format) already applied to rows 2-3 ---
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A4").Value = "2020-04-08"
$ws.Range("B4").Value = "Rendu de la version 1.0 de l'application Bataille Navale"

$ws.Rows.Item(4).RowHeight = 15.6

# --- Move the active selection down to B8 ---
$ws.Range("B8").Select() | Out-Null

# --- Page setup: portrait orientation, paper size 9 (A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "done"
